# B6-PowerPoint.pptx edit — Wed, May 06, 2020 4:08:07 PM
#
# 1) Three tables (on the "Activities" slides) had their table style
#    switched to a different built-in PowerPoint table style
#    ({EDCF7590-E278-4D9D-9AFD-E3D40A9DD64D} -> {62AFA160-25F6-4136-8DF4-26A162B7E33D}).
#
# 2) The presentation's theme colours were changed from the "Integral"
#    palette to the standard "Office" palette (the deck's notes master
#    already carried the Office theme; the main slide master is brought
#    in line with it here).

$p = $ppt.ActivePresentation

$oldStyleId = "{EDCF7590-E278-4D9D-9AFD-E3D40A9DD64D}"
$newStyleId = "{62AFA160-25F6-4136-8DF4-26A162B7E33D}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.StyleId -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# Re-point the deck's colour scheme at the standard Office palette.
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
